# Weekly fruit/vegetable price update.
# Inserts two new rows of this week's data above the existing history
# (rows 260-262 shift down to become rows 262-264, unchanged), then
# fills the two newly-opened rows (260-261) with the new entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: push the existing rows 260-262 down to 262-264.
$ws.Rows("260:261").Insert()

# New row 260: Choclero, Región Metropolitana
$ws.Range("A260").Value = 4
$ws.Range("B260").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C260").Value = "Los Lagos"
$ws.Range("D260").Value = 44656
$ws.Range("E260").Value = 10
$ws.Range("F260").Value = 100112024
$ws.Range("G260").Value = "Choclo"
$ws.Range("H260").Value = "Choclero"
$ws.Range("I260").Value = "Primera"
$ws.Range("J260").Value = 6000
$ws.Range("K260").Value = 350
$ws.Range("L260").Value = 350
$ws.Range("M260").Value = 350
$ws.Range("N260").Value = "$/unidad"
$ws.Range("O260").Value = "Región Metropolitana"
$ws.Range("P260").Value = 350
$ws.Range("Q260").Value = 1
$ws.Range("R260").Value = "Hortaliza"

# New row 261: Dulce o Americano, Región de Los Lagos
$ws.Range("A261").Value = 4
$ws.Range("B261").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C261").Value = "Los Lagos"
$ws.Range("D261").Value = 44656
$ws.Range("E261").Value = 10
$ws.Range("F261").Value = 100112024
$ws.Range("G261").Value = "Choclo"
$ws.Range("H261").Value = "Dulce o Americano"
$ws.Range("I261").Value = "Primera"
$ws.Range("J261").Value = 15000
$ws.Range("K261").Value = 180
$ws.Range("L261").Value = 200
$ws.Range("M261").Value = 190
$ws.Range("N261").Value = "$/unidad"
$ws.Range("O261").Value = "Región de Los Lagos"
$ws.Range("P261").Value = 190
$ws.Range("Q261").Value = 1
$ws.Range("R261").Value = "Hortaliza"
